$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Structural change: insert two banner rows on Sheet1 ---
# Insert a new row at row 2 (pushes old row2..46 down to 3..47)
$ws1.Rows("2:2").Insert() | Out-Null
# Insert a new row at row 21 (pushes old-row-now-at-21..47 down to 22..48)
$ws1.Rows("21:21").Insert() | Out-Null

# --- Re-point the hyperlink that moved from L39 to L41 because of the inserts ---
$ws1.Hyperlinks.Delete() | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("L41"), "https://forum.arduino.cc/index.php?topic=275431.0") | Out-Null

# --- New content appended at the bottom of Sheet1 (rows 50-53) ---
$ws1.Range("D50").Value2 = "Measure the time between pulses."
$ws1.Range("D51").Value2 = "As long as the incoming encoder is moving"
$ws1.Range("D52").Value2 = "Multiply the incoming pulses time delta by 7/3."
$ws1.Range("D53").Value2 = "Increment the output of the virtual encoder at that new rate."

# --- New notes added on Sheet2 near the blender-encoder wiring rows ---
$ws2.Range("B40").Value2 = "The blender encoders need the internal pull-up on the "
$ws2.Range("B41").Value2 = "Arduino to be enabled to drive a signal."

# --- New section-banner rows on Sheet1 (green fill, rgb FF92D050) ---
$bannerColor = 5296274  # RGB(146,208,80) == FF92D050

$band2 = $ws1.Range("A2:AK2")
$band2.Interior.Color = $bannerColor
$ws1.Range("A2").Value2 = "Input encoder triggers pin state-change interrupt for precision."

$band21 = $ws1.Range("A21:AK21")
$band21.Interior.Color = $bannerColor
$ws1.Range("A21").Value2 = "Output encoder is driven by time-based interrupt with 1ms period to allow output state change in between the input encoder pulses."

# --- View/selection changes ---
$ws2.Range("B42").Select() | Out-Null
$ws1.Range("C18").Select() | Out-Null
$ws1.Activate() | Out-Null
